$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H2").Value = 67528.53
$ws.Range("I2").Value = 342.25
$ws.Range("J2").Value = 144312.86
$ws.Range("K2").Value = 342.25
$ws.Range("L2").Value = 144312.86
$ws.Range("M2").Value = -229.25
$ws.Range("N2").Value = -144538.86

$ws.Range("H69").Value = 12800.5
$ws.Range("I69").Value = 12688.125
$ws.Range("J69").Value = 13250
$ws.Range("K69").Value = 38064.375
$ws.Range("L69").Value = 39750
$ws.Range("M69").Value = -37190.375
$ws.Range("N69").Value = -41498

$ws.Range("H72").Value = 12800.5
$ws.Range("I72").Value = 12688.125
$ws.Range("J72").Value = 13250
$ws.Range("K72").Value = 114193.125
$ws.Range("L72").Value = 119250
$ws.Range("M72").Value = -109825.125
$ws.Range("N72").Value = -127986

$ws.Range("H87").Value = 67372.875
$ws.Range("J87").Value = 67372.875
$ws.Range("L87").Value = 67372.875
$ws.Range("N87").Value = -69868.875

$ws.Range("H90").Value = 67372.875
$ws.Range("J90").Value = 67372.875
$ws.Range("L90").Value = 202118.625
$ws.Range("N90").Value = -214598.625

$ws.Range("H116").Value = 7280.4614
$ws.Range("I116").Value = 6780.7144
$ws.Range("K116").Value = 6780.7144
$ws.Range("M116").Value = -3338.7144

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H97").Value = 32259682
$ws.Range("I97").Value = 324.29413
$ws.Range("J97").Value = 71431760
$ws.Range("K97").Value = 324.29413
$ws.Range("L97").Value = 71431760
$ws.Range("M97").Value = 171.70587
$ws.Range("N97").Value = -71432752

$ws.Range("H110").Value = 55556252
$ws.Range("I110").Value = 71429224
$ws.Range("J110").Value = 849
$ws.Range("K110").Value = 71429224
$ws.Range("L110").Value = 849
$ws.Range("M110").Value = -71427179
$ws.Range("N110").Value = -4939

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H86").Value = 29413582
$ws.Range("J86").Value = 1923.4286
$ws.Range("L86").Value = 1923.4286
$ws.Range("N86").Value = -4169.4286

$ws.Range("H89").Value = 29413582
$ws.Range("J89").Value = 1923.4286
$ws.Range("L89").Value = 9617.143
$ws.Range("N89").Value = -20849.143

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 1883.64
$ws.Range("J31").Value = 2288.7273
$ws.Range("L31").Value = 2288.7273
$ws.Range("N31").Value = -2878.7273

$ws.Range("H34").Value = 1883.64
$ws.Range("J34").Value = 2288.7273
$ws.Range("L34").Value = 2288.7273
$ws.Range("N34").Value = -2692.7273

$ws.Range("H99").Value = 4745.125
$ws.Range("I99").Value = 4745.125
$ws.Range("K99").Value = 4745.125
$ws.Range("M99").Value = -3247.125

$ws.Range("H105").Value = 2144.3333
$ws.Range("I105").Value = 2144.3333
$ws.Range("J105").Value = 0
$ws.Range("K105").Value = 2144.3333
$ws.Range("L105").Value = 0
$ws.Range("M105").Value = -397.3332999999998
$ws.Range("N105").ClearContents()

$ws.Range("H107").Value = 660.3333
$ws.Range("I107").Value = 576.2941
$ws.Range("J107").Value = 770.2308
$ws.Range("K107").Value = 576.2941
$ws.Range("L107").Value = 770.2308
$ws.Range("M107").Value = 1343.7059
$ws.Range("N107").Value = -4610.2308

$ws.Range("H126").Value = 4745.125
$ws.Range("I126").Value = 4745.125
$ws.Range("K126").Value = 14235.375
$ws.Range("M126").Value = -11765.375

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 2487.625
$ws.Range("I5").Value = 925.5
$ws.Range("J5").Value = 3008.3333
$ws.Range("K5").Value = 2776.5
$ws.Range("L5").Value = 9024.999899999999
$ws.Range("M5").Value = -2664.5
$ws.Range("N5").Value = -9248.999899999999

$ws.Range("H74").Value = 10666.333
$ws.Range("I74").Value = 11999
$ws.Range("J74").Value = 10000
$ws.Range("K74").Value = 35997
$ws.Range("L74").Value = 30000
$ws.Range("M74").Value = -34936
$ws.Range("N74").Value = -32122

$ws.Range("H76").Value = 4333
$ws.Range("J76").Value = 7000
$ws.Range("L76").Value = 21000
$ws.Range("N76").Value = -21766

$ws.Range("H77").Value = 10666.333
$ws.Range("I77").Value = 11999
$ws.Range("J77").Value = 10000
$ws.Range("K77").Value = 107991
$ws.Range("L77").Value = 90000
$ws.Range("M77").Value = -102687
$ws.Range("N77").Value = -100608

$ws.Range("H79").Value = 4333
$ws.Range("J79").Value = 7000
$ws.Range("L79").Value = 21000
$ws.Range("N79").Value = -23652

$ws.Range("H80").Value = 14500.2
$ws.Range("I80").Value = 19667
$ws.Range("J80").Value = 6750
$ws.Range("K80").Value = 59001
$ws.Range("L80").Value = 20250
$ws.Range("M80").Value = -58065
$ws.Range("N80").Value = -22122

$ws.Range("H83").Value = 14500.2
$ws.Range("I83").Value = 19667
$ws.Range("J83").Value = 6750
$ws.Range("K83").Value = 177003
$ws.Range("L83").Value = 60750
$ws.Range("M83").Value = -172323
$ws.Range("N83").Value = -70110

$ws.Range("H87").Value = 5899.5
$ws.Range("I87").Value = 4299
$ws.Range("J87").Value = 7500
$ws.Range("K87").Value = 12897
$ws.Range("L87").Value = 22500
$ws.Range("M87").Value = -11649
$ws.Range("N87").Value = -24996

$ws.Range("H90").Value = 5899.5
$ws.Range("I90").Value = 4299
$ws.Range("J90").Value = 7500
$ws.Range("K90").Value = 38691
$ws.Range("L90").Value = 67500
$ws.Range("M90").Value = -32451
$ws.Range("N90").Value = -79980

$ws.Range("H121").Value = 10000900
$ws.Range("J121").Value = 14286926
$ws.Range("L121").Value = 42860778
$ws.Range("N121").Value = -42863398

$ws.Range("H135").Value = 2487.625
$ws.Range("I135").Value = 925.5
$ws.Range("J135").Value = 3008.3333
$ws.Range("K135").Value = 8329.5
$ws.Range("L135").Value = 27074.9997
$ws.Range("M135").Value = -5794.5
$ws.Range("N135").Value = -32144.9997

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 7750
$ws.Range("I80").Value = 7500
$ws.Range("J80").Value = 8000
$ws.Range("K80").Value = 7500
$ws.Range("L80").Value = 8000
$ws.Range("M80").Value = -6502
$ws.Range("N80").Value = -9996

$ws.Range("H83").Value = 7750
$ws.Range("I83").Value = 7500
$ws.Range("J83").Value = 8000
$ws.Range("K83").Value = 37500
$ws.Range("L83").Value = 40000
$ws.Range("M83").Value = -32508
$ws.Range("N83").Value = -49984

$ws.Range("H122").Value = 148236.08
$ws.Range("I122").Value = 205535.7
$ws.Range("J122").Value = 4987
$ws.Range("K122").Value = 616607.1000000001
$ws.Range("L122").Value = 14961
$ws.Range("M122").Value = -614157.1000000001
$ws.Range("N122").Value = -19861

$ws.Range("H132").Value = 7448.5747
$ws.Range("I132").Value = 7374.7837
$ws.Range("K132").Value = 22124.3511
$ws.Range("M132").Value = -19594.3511

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H39").Value = 30778
$ws.Range("I39").Value = 24890
$ws.Range("K39").Value = 24890
$ws.Range("M39").Value = -24430

$ws.Range("H40").Value = 4711.9375
$ws.Range("I40").Value = 5116.9165
$ws.Range("K40").Value = 5116.9165
$ws.Range("M40").Value = -4980.9165

$ws.Range("H61").Value = 58827452
$ws.Range("I61").Value = 62501040
$ws.Range("J61").Value = 49999
$ws.Range("K61").Value = 62501040
$ws.Range("L61").Value = 49999
$ws.Range("M61").Value = -62500838
$ws.Range("N61").Value = -50403

$ws.Range("H88").Value = 58632.145
$ws.Range("J88").Value = 63987.5
$ws.Range("L88").Value = 63987.5
$ws.Range("N88").Value = -64843.5

$ws.Range("H91").Value = 58632.145
$ws.Range("J91").Value = 63987.5
$ws.Range("L91").Value = 63987.5
$ws.Range("N91").Value = -66951.5

$ws.Range("H113").Value = 58827452
$ws.Range("I113").Value = 62501040
$ws.Range("J113").Value = 49999
$ws.Range("K113").Value = 62501040
$ws.Range("L113").Value = 49999
$ws.Range("M113").Value = -62498870
$ws.Range("N113").Value = -54339

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H39").Value = 9874.5
$ws.Range("J39").Value = 9999
$ws.Range("L39").Value = 9999
$ws.Range("N39").Value = -10825

$ws.Range("H56").Value = 34500
$ws.Range("I56").Value = 20000
$ws.Range("J56").Value = 39333.332
$ws.Range("K56").Value = 20000
$ws.Range("L56").Value = 39333.332
$ws.Range("N56").Value = -40761.332
$ws.Range("M56").Value = -19286
